# data model for program
# Rebuild the "program" sheet header/data row into the new wider schema,
# rename a couple of columns on "structures", and rename a couple of
# columns on "sections".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "program": new 14-column header schema + matching data row
# ---------------------------------------------------------------
$wsProgram = $wb.Worksheets.Item("program")

$programHeaders = @(
    "REPROG_ID_PRE",
    "REPROG_TITLE",
    "CED_ID_PRE",
    "CED_NAME_PRE",
    "REPROG_ACTIVE_IND",
    "REPROG_COMMENT",
    "REPROG_UW_DEPARTMENT_CD",
    "REPROG_UW_DEPARTMENT_NAME",
    "REPROG_UW_DEPARTMENT_LOB_CD",
    "REPROG_UW_DEPARTMENT_LOB_NAME",
    "BUSPAR_CED_REG_CLASS_CD",
    "BUSPAR_CED_REG_CLASS_NAME",
    "REPROG_MAIN_CURRENCY_CD",
    "REPROG_MANAGEMENT_REPORTING_LOB_CD"
)

# Keep the previous "program_name" value, now living in column B
# ("REPROG_TITLE"), and remember it before we overwrite A1.
$previousProgramName = $wsProgram.Cells.Item(2, 1).Value2

for ($i = 0; $i -lt $programHeaders.Count; $i++) {
    $wsProgram.Cells.Item(1, $i + 1).Value = $programHeaders[$i]
}

# Give every new header cell the same style as the pre-existing A1
# header cell (bold / bordered / centered).
$wsProgram.Range("A1").Copy()
$wsProgram.Range("A1:N1").PasteSpecial(-4122)

# Data row 2
$wsProgram.Cells.Item(2, 1).Value = 1
$wsProgram.Cells.Item(2, 2).Value = $previousProgramName
$wsProgram.Cells.Item(2, 5).Value = $true

$programEmptyCols = @(3, 4, 6, 7, 8, 9, 10, 11, 12, 13, 14)
foreach ($col in $programEmptyCols) {
    $wsProgram.Cells.Item(2, $col).Style = "Normal"
}

# ---------------------------------------------------------------
# Sheet "structures": rename columns + fix a typo in the data
# ---------------------------------------------------------------
$wsStructures = $wb.Worksheets.Item("structures")
$wsStructures.Range("B1").Value = "contract_order"
$wsStructures.Range("C1").Value = "type_of_participation"
$wsStructures.Range("C2").Value = "quota_share"

# ---------------------------------------------------------------
# Sheet "sections": rename columns
# ---------------------------------------------------------------
$wsSections = $wb.Worksheets.Item("sections")
$wsSections.Range("B1").Value = "cession_PCT"
$wsSections.Range("C1").Value = "attachment_point_100"
$wsSections.Range("D1").Value = "limit_occurrence_100"
